$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate row 9 (26-09-2025 entry with its values and formatting) into new row 10
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0
